$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.087075292901118928
$ws.Range("B1").Value = 0.086740928163777653
$ws.Range("A2").Value = -0.040628738363061245
$ws.Range("B2").Value = 0.039299331779030666
$ws.Range("A3").Value = 0.11821712792044181
$ws.Range("B3").Value = -0.11882502799748096
$ws.Range("A4").Value = -0.16517252606644206
$ws.Range("B4").Value = 0.16426490919332437
$ws.Range("A5").Value = -0.15826490930815318
$ws.Range("B5").Value = 0.15644552553691415
$ws.Range("A6").Value = -0.096847277183844671
$ws.Range("B6").Value = 0.096705582374324806
$ws.Range("A7").Value = -0.076705582518151516
$ws.Range("B7").Value = 0.076368209563423051
$ws.Range("A8").Value = -0.05636820970879608
$ws.Range("B8").Value = 0.056103435993489548
$ws.Range("A9").Value = -0.050103436116111233
$ws.Range("B9").Value = 0.049886643830347133
$ws.Range("A10").Value = -0.043886643954600402
$ws.Range("B10").Value = 0.0438559074501228
$ws.Range("A11").Value = -0.039355907571941628
$ws.Range("B11").Value = 0.039308028267377892
$ws.Range("A12").Value = -0.03330802839223912
$ws.Range("B12").Value = 0.033173052749885024
$ws.Range("A13").Value = -0.039154628369044708
$ws.Range("B13").Value = 0.039087049569597632
$ws.Range("A14").Value = -0.027087049707210653
$ws.Range("B14").Value = 0.027054115867723105
$ws.Range("A15").Value = -0.021054115995303491
$ws.Range("B15").Value = 0.021028164977879626
$ws.Range("A16").Value = -0.01502816510592897
$ws.Range("B16").Value = 0.015004778516943595
$ws.Range("A17").Value = -0.0090047786456155521
$ws.Range("B17").Value = 0.008999999865896946
$ws.Range("A18").Value = -0.03611218992131171
$ws.Range("B18").Value = 0.036097097023503721
$ws.Range("A19").Value = -0.027097097142444238
$ws.Range("B19").Value = 0.027014160412450305
$ws.Range("A20").Value = -0.018014160532429102
$ws.Range("B20").Value = 0.018004327213603233
$ws.Range("A21").Value = -0.0090043273337601093
$ws.Range("B21").Value = 0.0089999998797027914
$ws.Range("A22").Value = -0.14851645981891082
$ws.Range("B22").Value = 0.14766969563770438
$ws.Range("A23").Value = -0.13406432249458167
$ws.Range("B23").Value = 0.13270581630918699
$ws.Range("A24").Value = -0.042127364865700123
$ws.Range("B24").Value = 0.041999999816268563
$ws.Range("A25").Value = -0.053598248471804055
$ws.Range("B25").Value = 0.05353001685056924
$ws.Range("A26").Value = -0.047530016969549393
$ws.Range("B26").Value = 0.047449290839118419
$ws.Range("A27").Value = -0.041449290958476048
$ws.Range("B27").Value = 0.041197973175720826
$ws.Range("A28").Value = -0.035197973296395624
$ws.Range("B28").Value = 0.035040665610156729
$ws.Range("A29").Value = -0.023040665742163569
$ws.Range("B29").Value = 0.022983278673935459
$ws.Range("A30").Value = -0.0029832788202943661
$ws.Range("B30").Value = 0.0029394650221687968
$ws.Range("A31").Value = 0.012060534840147241
$ws.Range("B31").Value = -0.012085512828253542
$ws.Range("A32").Value = -0.0060010433419055431
$ws.Range("B32").Value = 0.0059999998781368191
